# Aggiornamento dati Fanano fino al 6 gennaio 2022
# Adds rows 465..491 (dates 44539..44566, serial 44549 absent as in source data)
# to the existing "Sheet1" sheet, replicating column A's date style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 464

$newData = @(
    @(44539,0,9,303.4389750505732),
    @(44540,1,7,236.0080917060014),
    @(44541,1,7,236.0080917060014),
    @(44542,0,7,236.0080917060014),
    @(44543,4,10,337.1544167228591),
    @(44544,1,7,236.0080917060014),
    @(44545,1,8,269.7235333782872),
    @(44546,3,11,370.8698583951449),
    @(44547,2,12,404.5853000674309),
    @(44548,4,15,505.7316250842886),
    @(44550,3,18,606.8779501011463),
    @(44551,2,16,539.4470667565745),
    @(44552,1,16,539.4470667565745),
    @(44553,2,17,573.1625084288604),
    @(44554,3,17,573.1625084288604),
    @(44555,0,15,505.7316250842886),
    @(44556,3,14,472.0161834120027),
    @(44557,4,15,505.7316250842886),
    @(44558,10,23,775.4551584625758),
    @(44559,10,32,1078.894133513149),
    @(44560,0,30,1011.463250168577),
    @(44561,1,28,944.0323668240054),
    @(44562,17,45,1517.194875252866),
    @(44563,2,44,1483.47943358058),
    @(44564,6,46,1550.910316925152),
    @(44565,4,40,1348.617666891436),
    @(44566,6,36,1213.755900202293)
)

$styleSourceCell = $ws.Range("A$lastExistingRow")

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $lastExistingRow + 1 + $i
    $values = $newData[$i]

    # Copy the date cell above (format + style) down, then overwrite with the
    # correct serial date value, so the new cell keeps style index/format
    # used throughout column A ("s=2" in the original sheet).
    $destA = $ws.Range("A$row")
    $styleSourceCell.Copy($destA)
    $destA.Value = $values[0]

    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
